$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

[object[,]]$row2 = New-Object "object[,]" 1,11
$row2[0,0] = 9.011401415084759
$row2[0,1] = 6.36283652711164
$row2[0,2] = 12.10265317743298
$row2[0,3] = 31.37987966924902
$row2[0,4] = 41.09093900627082
$row2[0,5] = 17.34030167517576
$row2[0,6] = 28.2865733793764
$row2[0,7] = 0
$row2[0,8] = 15.9397738775668
$row2[0,9] = 9.522728301295974
$row2[0,10] = 17.65933568960175
$ws.Range("C2:M2").Value = $row2

[object[,]]$row3 = New-Object "object[,]" 1,11
$row3[0,0] = 8.985308260655922
$row3[0,1] = 6.365330603851305
$row3[0,2] = 12.12534011658075
$row3[0,3] = 31.45515730929046
$row3[0,4] = 41.20870801798418
$row3[0,5] = 17.41481715928546
$row3[0,6] = 28.38611856401548
$row3[0,7] = 0
$row3[0,8] = 15.46882273521691
$row3[0,9] = 9.551196242636509
$row3[0,10] = 17.47274325961215
$ws.Range("C3:M3").Value = $row3

[object[,]]$row4 = New-Object "object[,]" 1,11
$row4[0,0] = 8.970581150613253
$row4[0,1] = 6.367257055817283
$row4[0,2] = 12.14123325594021
$row4[0,3] = 31.51210915043865
$row4[0,4] = 41.29781531524849
$row4[0,5] = 17.46459853887502
$row4[0,6] = 28.45532809354622
$row4[0,7] = 0
$row4[0,8] = 15.17416551669826
$row4[0,9] = 9.569833560874054
$row4[0,10] = 17.35983836787511
$ws.Range("C4:M4").Value = $row4

[object[,]]$row5 = New-Object "object[,]" 1,11
$row5[0,0] = 8.964908981736908
$row5[0,1] = 6.368141874688033
$row5[0,2] = 12.14820322794811
$row5[0,3] = 31.5380012966001
$row5[0,4] = 41.33832020622656
$row5[0,5] = 17.48589506502505
$row5[0,6] = 28.48555596711674
$row5[0,7] = 0
$row5[0,8] = 15.05288617286968
$row5[0,9] = 9.577719943569036
$row5[0,10] = 17.31429106895969
$ws.Range("C5:M5").Value = $row5

[object[,]]$row6 = New-Object "object[,]" 1,11
$row6[0,0] = 8.963987114137279
$row6[0,1] = 6.368294836885053
$row6[0,2] = 12.14939037636605
$row6[0,3] = 31.54246231292979
$row6[0,4] = 41.34529831586293
$row6[0,5] = 17.48949226898883
$row6[0,6] = 28.49069726978021
$row6[0,7] = 0
$row6[0,8] = 15.03268050685307
$row6[0,9] = 9.579047089554441
$row6[0,10] = 17.3067571707639
$ws.Range("C6:M6").Value = $row6

[object[,]]$row7 = New-Object "object[,]" 1,11
$row7[0,0] = 8.970503315807616
$row7[0,1] = 6.367268584264428
$row7[0,2] = 12.14132525813534
$row7[0,3] = 31.51244749462066
$row7[0,4] = 41.29834464364377
$row7[0,5] = 17.46488166473758
$row7[0,6] = 28.455727573219
$row7[0,7] = 0
$row7[0,8] = 15.17253453408292
$row7[0,9] = 9.569938738537484
$row7[0,10] = 17.35922217144488
$ws.Range("C7:M7").Value = $row7

[object[,]]$row8 = New-Object "object[,]" 1,11
$row8[0,0] = 9.002138253462254
$row8[0,1] = 6.363614741159207
$row8[0,2] = 12.110067998358
$row8[0,3] = 31.40359974454744
$row8[0,4] = 41.12804181860746
$row8[0,5] = 17.36515679989189
$row8[0,6] = 28.31921222721141
$row8[0,7] = 0
$row8[0,8] = 15.77863182410469
$row8[0,9] = 9.532303946687062
$row8[0,10] = 17.59467997990779
$ws.Range("C8:M8").Value = $row8

[object[,]]$row9 = New-Object "object[,]" 1,11
$row9[0,0] = 9.074264601079033
$row9[0,1] = 6.359566147255791
$row9[0,2] = 12.06436192804085
$row9[0,3] = 31.27589300158283
$row9[0,4] = 40.92861707449882
$row9[0,5] = 17.20168674489323
$row9[0,6] = 28.11609211478219
$row9[0,7] = 0
$row9[0,8] = 16.91649066430951
$row9[0,9] = 9.467673124965261
$row9[0,10] = 18.06764754939774
$ws.Range("C9:M9").Value = $row9

[object[,]]$row10 = New-Object "object[,]" 1,11
$row10[0,0] = 9.133152522039612
$row10[0,1] = 6.358468064446734
$row10[0,2] = 12.04029929205396
$row10[0,3] = 31.23506958936909
$row10[0,4] = 40.86567960778043
$row10[0,5] = 17.10131916079234
$row10[0,6] = 28.00676901218935
$row10[0,7] = 0
$row10[0,8] = 17.71294061239163
$row10[0,9] = 9.425756976498194
$row10[0,10] = 18.4193314602773
$ws.Range("C10:M10").Value = $row10

[object[,]]$row11 = New-Object "object[,]" 1,11
$row11[0,0] = 9.161163107602661
$row11[0,1] = 6.35837106706431
$row11[0,2] = 12.03142103218811
$row11[0,3] = 31.22813380375931
$row11[0,4] = 40.85547373623983
$row11[0,5] = 17.05998094421797
$row11[0,6] = 27.96580722609561
$row11[0,7] = 0
$row11[0,8] = 18.06511563762236
$row11[0,9] = 9.407892606863792
$row11[0,10] = 18.5796678951841
$ws.Range("C11:M11").Value = $row11

[object[,]]$row12 = New-Object "object[,]" 1,11
$row12[0,0] = 9.171940301827609
$row12[0,1] = 6.358391803433554
$row12[0,2] = 12.02835645215267
$row12[0,3] = 31.22718840326759
$row12[0,4] = 40.8542766157456
$row12[0,5] = 17.04495142405447
$row12[0,6] = 27.95156513776574
$row12[0,7] = 0
$row12[0,8] = 18.19690416880603
$row12[0,9] = 9.401300565313168
$row12[0,10] = 18.64038711762499
$ws.Range("C12:M12").Value = $row12

[object[,]]$row13 = New-Object "object[,]" 1,11
$row13[0,0] = 9.169611761586554
$row13[0,1] = 6.358384788201281
$row13[0,2] = 12.02900323617007
$row13[0,3] = 31.22731713888725
$row13[0,4] = 40.85441553803854
$row13[0,5] = 17.04816048012433
$row13[0,6] = 27.954575851595
$row13[0,7] = 0
$row13[0,8] = 18.16859279725939
$row13[0,9] = 9.402712597697795
$row13[0,10] = 18.62731074972242
$ws.Range("C13:M13").Value = $row13

[object[,]]$row14 = New-Object "object[,]" 1,11
$row14[0,0] = 9.162046375517773
$row14[0,1] = 6.358371623702329
$row14[0,2] = 12.03116294646146
$row14[0,3] = 31.22802229508233
$row14[0,4] = 40.85532169241684
$row14[0,5] = 17.05873192492377
$row14[0,6] = 27.96461003833866
$row14[0,7] = 0
$row14[0,8] = 18.07599009268956
$row14[0,9] = 9.407346813742103
$row14[0,10] = 18.58466350600629
$ws.Range("C14:M14").Value = $row14

[object[,]]$row15 = New-Object "object[,]" 1,11
$row15[0,0] = 9.157434353726455
$row15[0,1] = 6.358371031697477
$row15[0,2] = 12.03252456519011
$row15[0,3] = 31.22867335113222
$row15[0,4] = 40.8562246268642
$row15[0,5] = 17.06528864634615
$row15[0,6] = 27.97092179476007
$row15[0,7] = 0
$row15[0,8] = 18.01906025820637
$row15[0,9] = 9.41020790212942
$row15[0,10] = 18.55853986737793
$ws.Range("C15:M15").Value = $row15

[object[,]]$row16 = New-Object "object[,]" 1,11
$row16[0,0] = 9.131346088000761
$row16[0,1] = 6.358482467719517
$row16[0,2] = 12.04092112609609
$row16[0,3] = 31.23575770021296
$row16[0,4] = 40.86671896236524
$row16[0,5] = 17.10410790604475
$row16[0,6] = 28.00962314193176
$row16[0,7] = 0
$row16[0,8] = 17.68971109848711
$row16[0,9] = 9.426948650605596
$row16[0,10] = 18.40885613254892
$ws.Range("C16:M16").Value = $row16

[object[,]]$row17 = New-Object "object[,]" 1,11
$row17[0,0] = 9.115650929954011
$row17[0,1] = 6.358653630179689
$row17[0,2] = 12.04660185897152
$row17[0,3] = 31.24308987670913
$row17[0,4] = 40.87788982312599
$row17[0,5] = 17.12903089319623
$row17[0,6] = 28.03561729752167
$row17[0,7] = 0
$row17[0,8] = 17.48498466534398
$row17[0,9] = 9.437526618656666
$row17[0,10] = 18.31708654341063
$ws.Range("C17:M17").Value = $row17

[object[,]]$row18 = New-Object "object[,]" 1,11
$row18[0,0] = 9.106738874080536
$row18[0,1] = 6.358789990563659
$row18[0,2] = 12.05006391313371
$row18[0,3] = 31.24840191740053
$row18[0,4] = 40.88604833367329
$row18[0,5] = 17.1437723159271
$row18[0,6] = 28.05139345052368
$row18[0,7] = 0
$row18[0,8] = 17.36628367854222
$row18[0,9] = 9.443724079516192
$row18[0,10] = 18.26433920132479
$ws.Range("C18:M18").Value = $row18

[object[,]]$row19 = New-Object "object[,]" 1,11
$row19[0,0] = 9.103741387942668
$row19[0,1] = 6.358842684041602
$row19[0,2] = 12.05126953311469
$row19[0,3] = 31.2503882422066
$row19[0,4] = 40.88910776105757
$row19[0,5] = 17.14883320595947
$row19[0,6] = 28.05687642296612
$row19[0,7] = 0
$row19[0,8] = 17.32593449397492
$row19[0,9] = 9.445841899409315
$row19[0,10] = 18.24648751763115
$ws.Range("C19:M19").Value = $row19

[object[,]]$row20 = New-Object "object[,]" 1,11
$row20[0,0] = 9.117309804057577
$row20[0,1] = 6.35863148886662
$row20[0,2] = 12.0459769894089
$row20[0,3] = 31.24219598705038
$row20[0,4] = 40.87652113139282
$row20[0,5] = 17.12633571641759
$row20[0,6] = 28.0327647298581
$row20[0,7] = 0
$row20[0,8] = 17.50687710565747
$row20[0,9] = 9.436388852314483
$row20[0,10] = 18.32685218348283
$ws.Range("C20:M20").Value = $row20

[object[,]]$row21 = New-Object "object[,]" 1,11
$row21[0,0] = 9.16426393914287
$row21[0,1] = 6.358373934089764
$row21[0,2] = 12.03052051510959
$row21[0,3] = 31.22776949576528
$row21[0,4] = 40.85498300520622
$row21[0,5] = 17.05560986460849
$row21[0,6] = 27.96162824488048
$row21[0,7] = 0
$row21[0,8] = 18.10323328871153
$row21[0,9] = 9.405980944189219
$row21[0,10] = 18.59719032133925
$ws.Range("C21:M21").Value = $row21

[object[,]]$row22 = New-Object "object[,]" 1,11
$row22[0,0] = 9.195940470798382
$row22[0,1] = 6.358540378868464
$row22[0,2] = 12.02215240582225
$row22[0,3] = 31.22814180490321
$row22[0,4] = 40.85646171851532
$row22[0,5] = 17.01302719434299
$row22[0,6] = 27.92253802185873
$row22[0,7] = 0
$row22[0,8] = 18.4837755375209
$row22[0,9] = 9.387114761090071
$row22[0,10] = 18.77386961421879
$ws.Range("C22:M22").Value = $row22

[object[,]]$row23 = New-Object "object[,]" 1,11
$row23[0,0] = 9.178945480157431
$row23[0,1] = 6.358421049172745
$row23[0,2] = 12.02646000032778
$row23[0,3] = 31.22704412001354
$row23[0,4] = 40.85424402188038
$row23[0,5] = 17.03542017576172
$row23[0,6] = 27.94272140341725
$row23[0,7] = 0
$row23[0,8] = 18.28155081799085
$row23[0,9] = 9.397091926139332
$row23[0,10] = 18.67958827867737
$ws.Range("C23:M23").Value = $row23

[object[,]]$row24 = New-Object "object[,]" 1,11
$row24[0,0] = 9.116559479992011
$row24[0,1] = 6.358641380654236
$row24[0,2] = 12.04625888202913
$row24[0,3] = 31.24259669898467
$row24[0,4] = 40.87713450990616
$row24[0,5] = 17.12755291994723
$row24[0,6] = 28.0340517853085
$row24[0,7] = 0
$row24[0,8] = 17.49698264251195
$row24[0,9] = 9.436902874971445
$row24[0,10] = 18.32243709589703
$ws.Range("C24:M24").Value = $row24

[object[,]]$row25 = New-Object "object[,]" 1,11
$row25[0,0] = 9.053697514982726
$row25[0,1] = 6.360330189086931
$row25[0,2] = 12.0750559924718
$row25[0,3] = 31.30118059918249
$row25[0,4] = 40.96799184336101
$row25[0,5] = 17.24245658822897
$row25[0,6] = 28.16407015009917
$row25[0,7] = 0
$row25[0,8] = 16.61503853002941
$row25[0,9] = 9.484178045465445
$row25[0,10] = 17.93876803462964
$ws.Range("C25:M25").Value = $row25

